$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refreshed query data (Power Query "쿼리2" re-ran; totals for a few BJs went up) ---
$ws.Range("C3").Value = 608452
$ws.Range("C5").Value = 514683
$ws.Range("C7").Value = 340688

# --- Refresh timestamp column (새로고침시간) updated for every data row ---
$ws.Range("D2:D12").Value = 46014.918761006942

# --- Column D widened slightly (best-fit after new timestamp text) ---
$ws.Columns("D").ColumnWidth = 19.5

# --- Active cell moved to E13 (just below/right of the table after refresh) ---
$ws.Range("E13").Select()
